$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test user data row (row 2) with new registration values
$ws.Range("A2").Value = "test993"
$ws.Range("B2").Value = 23071337
$ws.Range("C2").Value = "narendra676"
$ws.Range("D2").Value = "Y`$k2&9aE"
